$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (bold font, border, center/top align) from row 15 pattern down to new rows 16-40
$ws.Range("A15:F15").Copy($ws.Range("A16:F40"))

# Update index column A for new rows (14..38)
for ($i = 16; $i -le 40; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# Set column B (always a value in rows 2-40)
$ws.Range("B2").Value = "NSE:AARTECH"
$ws.Range("B3").Value = "NSE:ABB"
$ws.Range("B4").Value = "NSE:ADROITINFO"
$ws.Range("B5").Value = "NSE:AJANTPHARM"
$ws.Range("B6").Value = "NSE:ALKYLAMINE"
$ws.Range("B7").Value = "NSE:APOLLO"
$ws.Range("B8").Value = "NSE:ASIANENE"
$ws.Range("B9").Value = "NSE:BANKBEES"
$ws.Range("B10").Value = "NSE:BANSWRAS"
$ws.Range("B11").Value = "NSE:CELLO"
$ws.Range("B12").Value = "NSE:CSLFINANCE"
$ws.Range("B13").Value = "NSE:EQUITASBNK"
$ws.Range("B14").Value = "NSE:ESSARSHPNG"
$ws.Range("B15").Value = "NSE:ESTER"
$ws.Range("B16").Value = "NSE:EXCELINDUS"
$ws.Range("B17").Value = "NSE:FUSION"
$ws.Range("B18").Value = "NSE:GOKUL"
$ws.Range("B19").Value = "NSE:HDFCQUAL"
$ws.Range("B20").Value = "NSE:INDOSTAR"
$ws.Range("B21").Value = "NSE:ITETF"
$ws.Range("B22").Value = "NSE:J&KBANK"
$ws.Range("B23").Value = "NSE:JOCIL"
$ws.Range("B24").Value = "NSE:KANANIIND"
$ws.Range("B25").Value = "NSE:KICL"
$ws.Range("B26").Value = "NSE:LAGNAM"
$ws.Range("B27").Value = "NSE:LIBAS"
$ws.Range("B28").Value = "NSE:LPDC"
$ws.Range("B29").Value = "NSE:MCL"
$ws.Range("B30").Value = "NSE:MOTOGENFIN"
$ws.Range("B31").Value = "NSE:NIACL"
$ws.Range("B32").Value = "NSE:NIPPOBATRY"
$ws.Range("B33").Value = "NSE:NPBET"
$ws.Range("B34").Value = "NSE:OSWALAGRO"
$ws.Range("B35").Value = "NSE:PEL"
$ws.Range("B36").Value = "NSE:POLYMED"
$ws.Range("B37").Value = "NSE:PREMEXPLN"
$ws.Range("B38").Value = "NSE:PSUBNKBEES"
$ws.Range("B39").Value = "NSE:RITCO"
$ws.Range("B40").Value = "NSE:SAKAR"

# Set column D values where needed (rows 2-3 only; rest remain empty)
$ws.Range("D2").Value = "NSE:HINDZINC"
$ws.Range("D3").Value = "NSE:MFSL"

# Set column E values where needed (rows 2-7)
$ws.Range("E2").Value = "NSE:CROMPTON"
$ws.Range("E3").Value = "NSE:CUMMINSIND"
$ws.Range("E4").Value = "NSE:EXIDEIND"
$ws.Range("E5").Value = "NSE:GLENMARK"
$ws.Range("E6").Value = "NSE:HEROMOTOCO"
$ws.Range("E7").Value = "NSE:JSWENERGY"

# Set column F values where needed (rows 2-9)
$ws.Range("F2").Value = "NSE:ABB"
$ws.Range("F3").Value = "NSE:AXISBANK"
$ws.Range("F4").Value = "NSE:BANKINDIA"
$ws.Range("F5").Value = "NSE:GAIL"
$ws.Range("F6").Value = "NSE:HDFCLIFE"
$ws.Range("F7").Value = "NSE:LICHSGFIN"
$ws.Range("F8").Value = "NSE:PEL"
$ws.Range("F9").Value = "NSE:PETRONET"
